$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(95).Style = "Neutral"
$ws.Range("A95").NumberFormat = "d-mmm"
$ws.Range("C95").NumberFormat = "0"
$ws.Range("H95").NumberFormat = "m/d/yy"
$ws.Range("I95").NumberFormat = "0"

$ws.Range("F95").Style = "Normal"
$ws.Range("J95").Style = "Normal"

$ws.Range("F95").ClearContents()
$ws.Range("I95").ClearContents()
$ws.Range("J95").ClearContents()

$ws.Rows(96).Style = "Good"
$ws.Range("A96").NumberFormat = "d-mmm"
$ws.Range("C96").NumberFormat = "0"
$ws.Range("H96").NumberFormat = "m/d/yy"
$ws.Range("I96").NumberFormat = "0"

$ws.Range("A96").Value = 44281
$ws.Range("B96").Value = 2987
$ws.Range("C96").Formula = "=(AVERAGE(B90:B96))"
$ws.Range("D96").Formula = "=AVERAGE(B83:B96)"
$ws.Range("E96").Formula = "=(E95-B96)"
$ws.Range("G96").Formula = "=E96/C96"
$ws.Range("H96").Formula = "=A96+G96"
$ws.Range("I96").Formula = "=E96/85"
$ws.Range("J96").Value = "daily rate to achieve June 20 target"
